# Remove the stray "Footer Placeholder 5" shape (ftr placeholder idx=11)
# that was duplicated/erroneously left on slide 27 ("Tweaks to Oct 7").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)

$target = $null
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Footer Placeholder 5") {
        $target = $sh
        break
    }
}

if ($target -ne $null) {
    # Cut() (rather than Delete()) removes the placeholder outright instead
    # of leaving an empty layout-inherited placeholder behind.
    $target.Cut()
}
